$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab name / workbook.xml sheet name)
$ws.Name = "SCD0024"

# 2. Update B2 cell value from "DGS-339" to "SCD0024-003"
$ws.Range("B2").Value = "SCD0024-003"

# 3. Update column B width (target stored OOXML width is 12.7109375; the
#    ColumnWidth property is quantized internally to steps of 1/6, so use
#    the value that rounds to the closest achievable stored width)
$ws.Columns.Item(2).ColumnWidth = 11.8333333

# 4. Update selection to B3
$ws.Range("B3").Select() | Out-Null
